$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for the Gas6-Mertk ligand-receptor pair table.
# Keys are cell addresses, values are the new numeric values from the refreshed dataset.
$newValues = @{
    "G2" = 20.030895
    "H2" = 60.092685
    "I2" = 0.1333691355055119
    "J2" = 0.1333691355055119
    "M2" = 11.01658666666667
    "N2" = 33.04976
    "O2" = 0.1837409300120545
    "P2" = 0.1837409300120545
    "Q2" = 220.6720907784
    "R2" = 1986.0488170056
    "S2" = 0.02450536899268648
    "T2" = 0.02450536899268648
    "G3" = 20.030895
    "H3" = 60.092685
    "I3" = 0.1333691355055119
    "J3" = 0.1333691355055119
    "O3" = 0.06509859443116503
    "P3" = 0.06509859443116503
    "Q3" = 78.18314046259499
    "R3" = 703.648264163355
    "S3" = 0.008682143261908413
    "T3" = 0.008682143261908413
    "G4" = 20.030895
    "H4" = 60.092685
    "I4" = 0.1333691355055119
    "J4" = 0.1333691355055119
    "M4" = 4.001997666666667
    "N4" = 12.005993
    "O4" = 0.06674760480978428
    "P4" = 0.06674760480978428
    "Q4" = 80.16359505124501
    "R4" = 721.4723554612051
    "S4" = 0.00890207035054448
    "T4" = 0.00890207035054448
    "G5" = 20.030895
    "H5" = 60.092685
    "I5" = 0.1333691355055119
    "J5" = 0.1333691355055119
    "M5" = 41.03546066666667
    "N5" = 123.106382
    "O5" = 0.6844128707469962
    "P5" = 0.6844128707469963
    "Q5" = 821.97700389063
    "R5" = 7397.79303501567
    "S5" = 0.09127955290037257
    "T5" = 0.09127955290037258
    "I6" = 0.6531407302146811
    "J6" = 0.653140730214681
    "M6" = 11.01658666666667
    "N6" = 33.04976
    "O6" = 0.1837409300120545
    "P6" = 0.1837409300120545
    "Q6" = 1080.684297477867
    "R6" = 9726.158677300798
    "S6" = 0.1200086851983979
    "T6" = 0.1200086851983979
    "I7" = 0.6531407302146811
    "J7" = 0.653140730214681
    "O7" = 0.06509859443116503
    "P7" = 0.06509859443116503
    "S7" = 0.0425185435027205
    "T7" = 0.0425185435027205
    "I8" = 0.6531407302146811
    "J8" = 0.653140730214681
    "M8" = 4.001997666666667
    "N8" = 12.005993
    "O8" = 0.06674760480978428
    "P8" = 0.06674760480978428
    "Q8" = 392.5804033290767
    "R8" = 3533.22362996169
    "S8" = 0.04359557934554347
    "T8" = 0.04359557934554346
    "I9" = 0.6531407302146811
    "J9" = 0.653140730214681
    "M9" = 41.03546066666667
    "N9" = 123.106382
    "O9" = 0.6844128707469962
    "P9" = 0.6844128707469963
    "Q9" = 4025.419063458006
    "R9" = 36228.77157112206
    "S9" = 0.4470179221680193
    "T9" = 0.4470179221680193
    "G10" = 3.916733333333333
    "H10" = 11.7502
    "I10" = 0.0260782825067122
    "J10" = 0.02607828250671219
    "M10" = 11.01658666666667
    "N10" = 33.04976
    "O10" = 0.1837409300120545
    "P10" = 0.1837409300120545
    "Q10" = 43.14903221688889
    "R10" = 388.341289952
    "S10" = 0.004791647880900391
    "T10" = 0.004791647880900391
    "G11" = 3.916733333333333
    "H11" = 11.7502
    "I11" = 0.0260782825067122
    "J11" = 0.02607828250671219
    "O11" = 0.06509859443116503
    "P11" = 0.06509859443116503
    "Q11" = 15.28751023628889
    "R11" = 137.5875921266
    "S11" = 0.001697659536365803
    "T11" = 0.001697659536365803
    "G12" = 3.916733333333333
    "H12" = 11.7502
    "I12" = 0.0260782825067122
    "J12" = 0.02607828250671219
    "M12" = 4.001997666666667
    "N12" = 12.005993
    "O12" = 0.06674760480978428
    "P12" = 0.06674760480978428
    "Q12" = 15.67475766095556
    "R12" = 141.0728189486
    "S12" = 0.001740662894875936
    "T12" = 0.001740662894875936
    "G13" = 3.916733333333333
    "H13" = 11.7502
    "I13" = 0.0260782825067122
    "J13" = 0.02607828250671219
    "M13" = 41.03546066666667
    "N13" = 123.106382
    "O13" = 0.6844128707469962
    "P13" = 0.6844128707469963
    "Q13" = 160.7249566418222
    "R13" = 1446.5246097764
    "S13" = 0.01784831219457007
    "T13" = 0.01784831219457007
    "G14" = 28.14764533333333
    "H14" = 84.442936
    "I14" = 0.1874118517730947
    "J14" = 0.1874118517730947
    "M14" = 11.01658666666667
    "N14" = 33.04976
    "O14" = 0.1837409300120545
    "P14" = 0.1837409300120545
    "Q14" = 310.0909742772622
    "R14" = 2790.81876849536
    "S14" = 0.03443522794006973
    "T14" = 0.03443522794006973
    "G15" = 28.14764533333333
    "H15" = 84.442936
    "I15" = 0.1874118517730947
    "J15" = 0.1874118517730947
    "O15" = 0.06509859443116503
    "P15" = 0.06509859443116503
    "Q15" = 109.8638532520542
    "R15" = 988.7746792684879
    "S15" = 0.01220024813017031
    "T15" = 0.01220024813017031
    "G16" = 28.14764533333333
    "H16" = 84.442936
    "I16" = 0.1874118517730947
    "J16" = 0.1874118517730947
    "M16" = 4.001997666666667
    "N16" = 12.005993
    "O16" = 0.06674760480978428
    "P16" = 0.06674760480978428
    "Q16" = 112.6468109461609
    "R16" = 1013.821298515448
    "S16" = 0.0125092922188204
    "T16" = 0.0125092922188204
    "G17" = 28.14764533333333
    "H17" = 84.442936
    "I17" = 0.1874118517730947
    "J17" = 0.1874118517730947
    "M17" = 41.03546066666667
    "N17" = 123.106382
    "O17" = 0.6844128707469962
    "P17" = 0.6844128707469963
    "Q17" = 1155.051592935283
    "R17" = 10395.46433641755
    "S17" = 0.1282670834840343
    "T17" = 0.1282670834840343
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}
